# Auto-generated edit script: updates cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force target Price cells to Text format so values remain strings (matches source data)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "61.485.81"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").Value = "2.992.69"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "547.21"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "131.12"
$ws.Range("E6").Value = "  -5.83%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "2.987.15"
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Value = "6.00"
$ws.Range("E10").Value = "  -6.64%  "
$ws.Range("D11").Value = "0.145"
$ws.Range("E11").Value = "  -7.89%  "
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  -2.99%  "
$ws.Range("D13").Value = "34.00"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Value = "0.0000218"
$ws.Range("E14").Value = "  -3.10%  "
$ws.Range("D15").Value = "3.473.51"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "61.586.92"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").Value = "0.109"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "2.989.73"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").Value = "6.61"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "478.97"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "13.12"
$ws.Range("D22").Value = "0.663"
$ws.Range("E22").Value = "  -5.37%  "
$ws.Range("D23").Value = "6.97"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").Value = "80.79"
$ws.Range("E24").Value = "  +2.68%  "
$ws.Range("D25").Value = "12.01"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "2.71"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "1.91"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").Value = "25.46"
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("D33").Value = "2.31"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").Value = "54.81"
$ws.Range("E35").Value = "  -7.81%  "
$ws.Range("D36").Value = "5.85"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").Value = "446.81"
$ws.Range("E37").Value = "  -8.83%  "
$ws.Range("D38").Value = "3.126.09"
$ws.Range("E38").Value = "  -4.36%  "
$ws.Range("D39").Value = "0.0792"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "0.0381"
$ws.Range("E40").Value = "  -5.74%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").Value = "8.07"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D44").Value = "2.36"
$ws.Range("E44").Value = "  -9.48%  "
$ws.Range("D45").Value = "25.52"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "0.242"
$ws.Range("E46").Value = "  -4.53%  "
$ws.Range("D47").Value = "0.108"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "1.95"
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "114.50"
$ws.Range("E49").Value = "  -7.67%  "
$ws.Range("B50").Value = "BitgetToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D50").Value = "1.30"
$ws.Range("E50").Value = "  +9.05%  "
$ws.Range("D51").Value = "0.0₃0483"
$ws.Range("E51").Value = "  -8.89%  "
